$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Rarres2"
$ws.Range("C2").Value = "Cmklr1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.9064200000000001
$ws.Range("H2").Value = 2.71926
$ws.Range("I2").Value = 0.01741933661286065
$ws.Range("J2").Value = 0.01741933661286065
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.576447
$ws.Range("N2").Value = 13.729341
$ws.Range("O2").Value = 0.06514895840416586
$ws.Range("P2").Value = 0.06514895840416586
$ws.Range("Q2").Value = 4.148183089740001
$ws.Range("R2").Value = 37.33364780766
$ws.Range("S2").Value = 0.001134851636419422
$ws.Range("T2").Value = 0.001134851636419422

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Rarres2"
$ws.Range("C3").Value = "Cmklr1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.9064200000000001
$ws.Range("H3").Value = 2.71926
$ws.Range("I3").Value = 0.01741933661286065
$ws.Range("J3").Value = 0.01741933661286065
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 25.39955333333333
$ws.Range("N3").Value = 76.19865999999999
$ws.Range("O3").Value = 0.361580598136005
$ws.Range("P3").Value = 0.361580598136005
$ws.Range("Q3").Value = 23.0226631324
$ws.Range("R3").Value = 207.2039681916
$ws.Range("S3").Value = 0.006298494151610565
$ws.Range("T3").Value = 0.006298494151610565

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Rarres2"
$ws.Range("C4").Value = "Cmklr1"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.9064200000000001
$ws.Range("H4").Value = 2.71926
$ws.Range("I4").Value = 0.01741933661286065
$ws.Range("J4").Value = 0.01741933661286065
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 38.144755
$ws.Range("N4").Value = 114.434265
$ws.Range("O4").Value = 0.5430175541926078
$ws.Range("P4").Value = 0.5430175541926079
$ws.Range("Q4").Value = 34.5751688271
$ws.Range("R4").Value = 311.1765194439
$ws.Range("S4").Value = 0.009459005563173336
$ws.Range("T4").Value = 0.009459005563173337

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Rarres2"
$ws.Range("C5").Value = "Cmklr1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.9064200000000001
$ws.Range("H5").Value = 2.71926
$ws.Range("I5").Value = 0.01741933661286065
$ws.Range("J5").Value = 0.01741933661286065
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.125141333333333
$ws.Range("N5").Value = 6.375424
$ws.Range("O5").Value = 0.03025288926722125
$ws.Range("P5").Value = 0.03025288926722126
$ws.Range("Q5").Value = 1.92627060736
$ws.Range("R5").Value = 17.33643546624
$ws.Range("S5").Value = 0.0005269852616573261
$ws.Range("T5").Value = 0.0005269852616573263

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Rarres2"
$ws.Range("C6").Value = "Cmklr1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 42.58841866666666
$ws.Range("H6").Value = 127.765256
$ws.Range("I6").Value = 0.8184528149909583
$ws.Range("J6").Value = 0.8184528149909585
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.576447
$ws.Range("N6").Value = 13.729341
$ws.Range("O6").Value = 0.06514895840416586
$ws.Range("P6").Value = 0.06514895840416586
$ws.Range("Q6").Value = 194.9036408418106
$ws.Range("R6").Value = 1754.132767576296
$ws.Range("S6").Value = 0.0533213483996184
$ws.Range("T6").Value = 0.05332134839961841

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Rarres2"
$ws.Range("C7").Value = "Cmklr1"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 42.58841866666666
$ws.Range("H7").Value = 127.765256
$ws.Range("I7").Value = 0.8184528149909583
$ws.Range("J7").Value = 0.8184528149909585
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 25.39955333333333
$ws.Range("N7").Value = 76.19865999999999
$ws.Range("O7").Value = 0.361580598136005
$ws.Range("P7").Value = 0.361580598136005
$ws.Range("Q7").Value = 1081.726811306329
$ws.Range("R7").Value = 9735.541301756959
$ws.Range("S7").Value = 0.2959366583905277
$ws.Range("T7").Value = 0.2959366583905277

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Rarres2"
$ws.Range("C8").Value = "Cmklr1"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 42.58841866666666
$ws.Range("H8").Value = 127.765256
$ws.Range("I8").Value = 0.8184528149909583
$ws.Range("J8").Value = 0.8184528149909585
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 38.144755
$ws.Range("N8").Value = 114.434265
$ws.Range("O8").Value = 0.5430175541926078
$ws.Range("P8").Value = 0.5430175541926079
$ws.Range("Q8").Value = 1624.524795877426
$ws.Range("R8").Value = 14620.72316289684
$ws.Range("S8").Value = 0.4444342458184451
$ws.Range("T8").Value = 0.4444342458184453

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Rarres2"
$ws.Range("C9").Value = "Cmklr1"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 42.58841866666666
$ws.Range("H9").Value = 127.765256
$ws.Range("I9").Value = 0.8184528149909583
$ws.Range("J9").Value = 0.8184528149909585
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.125141333333333
$ws.Range("N9").Value = 6.375424
$ws.Range("O9").Value = 0.03025288926722125
$ws.Range("P9").Value = 0.03025288926722126
$ws.Range("Q9").Value = 90.50640882983821
$ws.Range("R9").Value = 814.5576794685439
$ws.Range("S9").Value = 0.02476056238236698
$ws.Range("T9").Value = 0.02476056238236699

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Rarres2"
$ws.Range("C10").Value = "Cmklr1"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 8.540438
$ws.Range("H10").Value = 25.621314
$ws.Range("I10").Value = 0.164127848396181
$ws.Range("J10").Value = 0.164127848396181
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 4.576447
$ws.Range("N10").Value = 13.729341
$ws.Range("O10").Value = 0.06514895840416586
$ws.Range("P10").Value = 0.06514895840416586
$ws.Range("Q10").Value = 39.084861863786
$ws.Range("R10").Value = 351.763756774074
$ws.Range("S10").Value = 0.01069275836812804
$ws.Range("T10").Value = 0.01069275836812804

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Rarres2"
$ws.Range("C11").Value = "Cmklr1"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 8.540438
$ws.Range("H11").Value = 25.621314
$ws.Range("I11").Value = 0.164127848396181
$ws.Range("J11").Value = 0.164127848396181
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 25.39955333333333
$ws.Range("N11").Value = 76.19865999999999
$ws.Range("O11").Value = 0.361580598136005
$ws.Range("P11").Value = 0.361580598136005
$ws.Range("Q11").Value = 216.9233104710266
$ws.Range("R11").Value = 1952.30979423924
$ws.Range("S11").Value = 0.05934544559386667
$ws.Range("T11").Value = 0.05934544559386667

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Rarres2"
$ws.Range("C12").Value = "Cmklr1"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 8.540438
$ws.Range("H12").Value = 25.621314
$ws.Range("I12").Value = 0.164127848396181
$ws.Range("J12").Value = 0.164127848396181
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 38.144755
$ws.Range("N12").Value = 114.434265
$ws.Range("O12").Value = 0.5430175541926078
$ws.Range("P12").Value = 0.5430175541926079
$ws.Range("Q12").Value = 325.77291510269
$ws.Range("R12").Value = 2931.95623592421
$ws.Range("S12").Value = 0.08912430281098933
$ws.Range("T12").Value = 0.08912430281098935

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Rarres2"
$ws.Range("C13").Value = "Cmklr1"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 8.540438
$ws.Range("H13").Value = 25.621314
$ws.Range("I13").Value = 0.164127848396181
$ws.Range("J13").Value = 0.164127848396181
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 2.125141333333333
$ws.Range("N13").Value = 6.375424
$ws.Range("O13").Value = 0.03025288926722125
$ws.Range("P13").Value = 0.03025288926722126
$ws.Range("Q13").Value = 18.14963779857067
$ws.Range("R13").Value = 163.346740187136
$ws.Range("S13").Value = 0.004965341623196941
$ws.Range("T13").Value = 0.004965341623196942
